# Updates "想去人数" (people interested) counts in column F across the
# workbook's sheets, reflecting newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F10").Value = 1141
$ws.Range("F11").Value = 668
$ws.Range("F12").Value = 449
$ws.Range("F13").Value = 741
$ws.Range("F15").Value = 215
$ws.Range("F17").Value = 233
$ws.Range("F19").Value = 279
$ws.Range("F20").Value = 1471
$ws.Range("F25").Value = 2209
$ws.Range("F27").Value = 775
$ws.Range("F30").Value = 45

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 437

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 28
$ws.Range("F5").Value = 1767
$ws.Range("F6").Value = 1769
$ws.Range("F7").Value = 595
$ws.Range("F8").Value = 579
$ws.Range("F9").Value = 447

# Sheet "全部类型" (All Types) - aggregated view of the above sheets
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1767
$ws.Range("F6").Value = 1769
$ws.Range("F7").Value = 595
$ws.Range("F12").Value = 579
$ws.Range("F14").Value = 447
$ws.Range("F19").Value = 1141
$ws.Range("F20").Value = 668
$ws.Range("F21").Value = 449
$ws.Range("F24").Value = 741
$ws.Range("F26").Value = 215
$ws.Range("F31").Value = 233
$ws.Range("F33").Value = 279
$ws.Range("F35").Value = 1471
$ws.Range("F37").Value = 437
$ws.Range("F41").Value = 2209
$ws.Range("F44").Value = 775
$ws.Range("F48").Value = 45
